$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value for every data row (2-219).
# All of them change from 45172 (2023-09-03) to 45175 (2023-09-06).
$ws.Range("C2:C219").Value = 45175
